$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3072.5173
$ws.Range("I137").Value = 1477.7778
$ws.Range("K137").Value = 4433.3334
$ws.Range("M137").Value = -1883.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2258.4119
$ws.Range("I2").Value = 1813
$ws.Range("J2").Value = 2822.6
$ws.Range("K2").Value = 1813
$ws.Range("L2").Value = 2822.6
$ws.Range("M2").Value = -1700
$ws.Range("N2").Value = -3048.6
$ws.Range("H32").Value = 8341964
$ws.Range("I32").Value = 9620659
$ws.Range("J32").Value = 30444.25
$ws.Range("K32").Value = 9620659
$ws.Range("L32").Value = 30444.25
$ws.Range("M32").Value = -9620372
$ws.Range("N32").Value = -31018.25
$ws.Range("H45").Value = 5499.125
$ws.Range("I45").Value = 6999
$ws.Range("K45").Value = 6999
$ws.Range("M45").Value = -6622
$ws.Range("H74").Value = 4034937.8
$ws.Range("I74").Value = 4809432.5
$ws.Range("K74").Value = 4809432.5
$ws.Range("M74").Value = -4808558.5
$ws.Range("H77").Value = 4034937.8
$ws.Range("I77").Value = 4809432.5
$ws.Range("K77").Value = 24047162.5
$ws.Range("M77").Value = -24042794.5
$ws.Range("H116").Value = 2258.4119
$ws.Range("I116").Value = 1813
$ws.Range("J116").Value = 2822.6
$ws.Range("K116").Value = 1813
$ws.Range("L116").Value = 2822.6
$ws.Range("M116").Value = 481
$ws.Range("N116").Value = -7410.6
$ws.Range("H122").Value = 2742.625
$ws.Range("I122").Value = 2848.7856
$ws.Range("K122").Value = 8546.356800000001
$ws.Range("M122").Value = -6096.356800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2258.4119
$ws.Range("I3").Value = 1813
$ws.Range("J3").Value = 2822.6
$ws.Range("K3").Value = 1813
$ws.Range("L3").Value = 2822.6
$ws.Range("M3").Value = -1699
$ws.Range("N3").Value = -3050.6
$ws.Range("H22").Value = 748
$ws.Range("I22").Value = 748
$ws.Range("K22").Value = 748
$ws.Range("M22").Value = -575
$ws.Range("H80").Value = 15059.8
$ws.Range("I80").Value = 25845
$ws.Range("K80").Value = 25845
$ws.Range("M80").Value = -24847
$ws.Range("H83").Value = 15059.8
$ws.Range("I83").Value = 25845
$ws.Range("K83").Value = 129225
$ws.Range("M83").Value = -124233
$ws.Range("H107").Value = 500
$ws.Range("I107").Value = 500
$ws.Range("K107").Value = 500
$ws.Range("M107").Value = 1420
$ws.Range("H134").Value = 481411.47
$ws.Range("I134").Value = 772967.4
$ws.Range("J134").Value = 5715.0527
$ws.Range("K134").Value = 2318902.2
$ws.Range("L134").Value = 17145.1581
$ws.Range("M134").Value = -2316367.2
$ws.Range("N134").Value = -22215.1581

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1417.8182
$ws.Range("I16").Value = 1399.5555
$ws.Range("K16").Value = 1399.5555
$ws.Range("M16").Value = -1112.5555
$ws.Range("H31").Value = 13764.814
$ws.Range("I31").Value = 4654.6665
$ws.Range("K31").Value = 4654.6665
$ws.Range("M31").Value = -4359.6665
$ws.Range("H34").Value = 13764.814
$ws.Range("I34").Value = 4654.6665
$ws.Range("K34").Value = 4654.6665
$ws.Range("M34").Value = -4452.6665
$ws.Range("H53").Value = 24950
$ws.Range("J53").Value = 24950
$ws.Range("L53").Value = 24950
$ws.Range("N53").Value = -26164
$ws.Range("H58").Value = 777261.4399999999
$ws.Range("I58").Value = 1034473.56
$ws.Range("J58").Value = 5625
$ws.Range("K58").Value = 1034473.56
$ws.Range("L58").Value = 5625
$ws.Range("M58").Value = -1034270.56
$ws.Range("N58").Value = -6031
$ws.Range("H69").Value = 25942.143
$ws.Range("I69").Value = 5653.8184
$ws.Range("K69").Value = 5653.8184
$ws.Range("M69").Value = -4904.8184
$ws.Range("H72").Value = 25942.143
$ws.Range("I72").Value = 5653.8184
$ws.Range("K72").Value = 16961.4552
$ws.Range("M72").Value = -13217.4552
$ws.Range("H113").Value = 1417.8182
$ws.Range("I113").Value = 1399.5555
$ws.Range("K113").Value = 1399.5555
$ws.Range("M113").Value = 770.4445000000001
$ws.Range("H136").Value = 777261.4399999999
$ws.Range("I136").Value = 1034473.56
$ws.Range("J136").Value = 5625
$ws.Range("K136").Value = 3103420.68
$ws.Range("L136").Value = 16875
$ws.Range("M136").Value = -3100870.68
$ws.Range("N136").Value = -21975

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 6662.6665
$ws.Range("J62").Value = 6662.6665
$ws.Range("L62").Value = 19987.9995
$ws.Range("N62").Value = -21359.9995
$ws.Range("H65").Value = 6662.6665
$ws.Range("J65").Value = 6662.6665
$ws.Range("L65").Value = 59963.9985
$ws.Range("N65").Value = -66827.9985
$ws.Range("H75").Value = 7260.6924
$ws.Range("I75").Value = 4116.5
$ws.Range("J75").Value = 7832.364
$ws.Range("K75").Value = 12349.5
$ws.Range("L75").Value = 23497.092
$ws.Range("M75").Value = -11351.5
$ws.Range("N75").Value = -25493.092
$ws.Range("H78").Value = 7260.6924
$ws.Range("I78").Value = 4116.5
$ws.Range("J78").Value = 7832.364
$ws.Range("K78").Value = 37048.5
$ws.Range("L78").Value = 70491.276
$ws.Range("M78").Value = -32056.5
$ws.Range("N78").Value = -80475.276
$ws.Range("H129").Value = 3043.3125
$ws.Range("I129").Value = 742.25
$ws.Range("J129").Value = 3810.3333
$ws.Range("K129").Value = 2226.75
$ws.Range("L129").Value = 11430.9999
$ws.Range("M129").Value = 2773.25
$ws.Range("N129").Value = -21430.9999
$ws.Range("H132").Value = 3600.6843
$ws.Range("I132").Value = 2785.6365
$ws.Range("J132").Value = 4721.375
$ws.Range("K132").Value = 25070.7285
$ws.Range("L132").Value = 42492.375
$ws.Range("M132").Value = -22540.7285
$ws.Range("N132").Value = -47552.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 19016.25
$ws.Range("J32").Value = 19016.25
$ws.Range("L32").Value = 19016.25
$ws.Range("N32").Value = -19608.25
$ws.Range("H139").Value = 138499.4
$ws.Range("J139").Value = 143888.22
$ws.Range("L139").Value = 143888.22
$ws.Range("N139").Value = -154168.22

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 42673.56
$ws.Range("I22").Value = 144007.86
$ws.Range("J22").Value = 3265.7778
$ws.Range("K22").Value = 144007.86
$ws.Range("L22").Value = 3265.7778
$ws.Range("M22").Value = -143712.86
$ws.Range("N22").Value = -3855.7778
$ws.Range("H27").Value = 42673.56
$ws.Range("I27").Value = 144007.86
$ws.Range("J27").Value = 3265.7778
$ws.Range("K27").Value = 144007.86
$ws.Range("L27").Value = 3265.7778
$ws.Range("M27").Value = -143900.86
$ws.Range("N27").Value = -3479.7778
$ws.Range("H46").Value = 3673.7334
$ws.Range("I46").Value = 1331.6666
$ws.Range("K46").Value = 1331.6666
$ws.Range("M46").Value = -1143.6666
$ws.Range("H55").Value = 1571.9333
$ws.Range("I55").Value = 963
$ws.Range("K55").Value = 963
$ws.Range("M55").Value = -790
$ws.Range("H93").Value = 2788.9565
$ws.Range("I93").Value = 1270.8422
$ws.Range("K93").Value = 1270.8422
$ws.Range("M93").Value = -22.84220000000005
$ws.Range("H122").Value = 78236.86
$ws.Range("I122").Value = 4052.2
$ws.Range("K122").Value = 12156.6
$ws.Range("M122").Value = -9706.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 5666.6665
$ws.Range("I100").Value = 7500
$ws.Range("K100").Value = 15000
$ws.Range("M100").Value = -14459
